# Rename the worksheet from "Property1" to "DataNode" to unify the
# DataNode/DataTable/Entity naming convention described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Adjust the header row height (row 1) and the description row height
# (row 8) to match the re-tuned layout.
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 81

# Move/restore the active selection to H13, matching where the editor
# left the cursor before saving.
$ws.Range("H13").Select() | Out-Null
